$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B -> shifts RawActivations/PercActivations/totalActivation right
$ws.Range("B1").EntireColumn.Insert()

# New header for the inserted column (segment names header) - copy header formatting from C1
$ws.Range("C1").Copy()
$ws.Range("B1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("B1").Value = "segments"

# The column-insert carried column A's formatting into the new B2:B20 cells - reset to default
$ws.Range("B2:B20").Style = "Normal"

# Move the segment-name labels from column A to column B, and put a 0-based index in column A
for ($r = 2; $r -le 20; $r++) {
    $name = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($r, 2).Value = $name
    $ws.Cells.Item($r, 1).Value = $r - 2
}

$wb.Save()
